# Commit: "adding pulse ox profile"
#
# The "Pulse oximetry" / "Inhaled oxygen concentration" rows in the USCDI
# table on Sheet1 were pointing at the generic FHIR-core "[Oxygen
# saturation]" profile (or left blank). This adds a dedicated
# "[US Core Pulse Oximetry Profile]" entry and wires it up to both rows,
# and fills in the previously-blank FHIR Resource cell for the
# "Inhaled oxygen concentration" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 54 = "  Pulse oximetry" (B54/C54 carry no special cell formatting)
$ws.Range("B54").Value = "[US Core Pulse Oximetry Profile]"
$ws.Range("C54").Value = "Observation"

# Row 55 = "  Inhaled oxygen concentration" -- B55/C55 use the
# quote-prefixed "forced text" style (they used to hold the literal "-").
# Preserve that formatting across the value change by round-tripping the
# cell formats through a scratch cell instead of a plain .Value= write
# (which would otherwise drop the cell style).
$ws.Range("B55").Copy() | Out-Null
$ws.Range("ZZ1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B55").Value = "[US Core Pulse Oximetry Profile]"
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("B55").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C55").Copy() | Out-Null
$ws.Range("ZZ1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C55").Value = "Observation"
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("C55").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Clean up the scratch cell used to stage the copied format.
$ws.Range("ZZ1").Clear() | Out-Null
